# Trade #15 closed/recorded at 2026-02-16 22:53:31 - base_strategy UP +0.000%
#
# Appends a new trade-log row (row 16) to both the "All Trades" sheet and the
# "base_strategy" sheet, mirroring the existing row layout:
#   A: Trade #            I: P&L %
#   B: Date               J: P&L $
#   C: Time               K: Capital After
#   D: Strategy           L: Entry Slippage (bps)
#   E: Side               M: Exit Slippage (bps)
#   F: Entry Price        N: Confidence
#   G: Exit Price         O: Entry Reason
#   H: Status             P: Exit Reason
#                         Q: Duration (min)

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 16

    $ws.Cells.Item($row, 1).Value = 15              # A16 Trade #

    # B16 Date - force text so "2026-02-16" is stored as a literal string
    # rather than being auto-parsed into a date serial number.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"

    $ws.Cells.Item($row, 3).Value = "22:53:31"       # C16 Time
    $ws.Cells.Item($row, 4).Value = "base_strategy"  # D16 Strategy
    $ws.Cells.Item($row, 5).Value = "UP"             # E16 Side
    $ws.Cells.Item($row, 6).Value = 49.999998        # F16 Entry Price
    # G16 Exit Price intentionally left blank (trade still OPEN)
    $ws.Cells.Item($row, 8).Value = "OPEN"           # H16 Status
    $ws.Cells.Item($row, 9).Value = 0                # I16 P&L %
    $ws.Cells.Item($row, 10).Value = 0               # J16 P&L $
    $ws.Cells.Item($row, 11).Value = 100             # K16 Capital After
    $ws.Cells.Item($row, 12).Value = 0               # L16 Entry Slippage (bps)
    $ws.Cells.Item($row, 13).Value = 0               # M16 Exit Slippage (bps)
    $ws.Cells.Item($row, 14).Value = 0.6             # N16 Confidence
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"  # O16 Entry Reason
    # P16 Exit Reason intentionally left blank (trade still OPEN)
    $ws.Cells.Item($row, 17).Value = 0               # Q16 Duration (min)
}
